$wb = $excel.ActiveWorkbook

# Sheet 1: 展览 (Exhibitions) - update "想去人数" (F) / "最低票价" (G) counts
$ws = $wb.Worksheets.Item(1)
$ws.Cells.Item(4, 6).Value = 553
$ws.Cells.Item(6, 6).Value = 1594
$ws.Cells.Item(8, 6).Value = 39
$ws.Cells.Item(9, 6).Value = 723
$ws.Cells.Item(10, 6).Value = 2664
$ws.Cells.Item(11, 6).Value = 2664
$ws.Cells.Item(13, 6).Value = 1741
$ws.Cells.Item(14, 6).Value = 604
$ws.Cells.Item(15, 6).Value = 267
$ws.Cells.Item(15, 7).Value = 139
$ws.Cells.Item(17, 6).Value = 4974
$ws.Cells.Item(18, 6).Value = 179
$ws.Cells.Item(19, 6).Value = 63
$ws.Cells.Item(21, 6).Value = 3382
$ws.Cells.Item(22, 6).Value = 857
$ws.Cells.Item(24, 6).Value = 68
$ws.Cells.Item(25, 6).Value = 36
$ws.Cells.Item(26, 6).Value = 2408
$ws.Cells.Item(27, 6).Value = 56
$ws.Cells.Item(28, 6).Value = 363
$ws.Cells.Item(32, 6).Value = 1289
$ws.Cells.Item(33, 6).Value = 800
$ws.Cells.Item(34, 6).Value = 3
$ws.Cells.Item(35, 6).Value = 55
$ws.Cells.Item(38, 6).Value = 1410
$ws.Cells.Item(39, 6).Value = 2
$ws.Cells.Item(40, 6).Value = 1372
# Sheet 2: 演出 (Performances)
$ws = $wb.Worksheets.Item(2)
$ws.Cells.Item(8, 6).Value = 18
$ws.Cells.Item(9, 6).Value = 114
$ws.Cells.Item(10, 6).Value = 218
$ws.Cells.Item(15, 6).Value = 25
$ws.Cells.Item(17, 6).Value = 324
$ws.Cells.Item(19, 6).Value = 515
# Sheet 3: 本地生活 (Local life)
$ws = $wb.Worksheets.Item(3)
$ws.Cells.Item(3, 6).Value = 856
$ws.Cells.Item(4, 6).Value = 242
$ws.Cells.Item(6, 6).Value = 20
$ws.Cells.Item(7, 6).Value = 27
# Sheet 4: 全部类型 (All types, aggregated)
$ws = $wb.Worksheets.Item(4)
$ws.Cells.Item(6, 6).Value = 856
$ws.Cells.Item(7, 6).Value = 242
$ws.Cells.Item(10, 6).Value = 553
$ws.Cells.Item(11, 6).Value = 20
$ws.Cells.Item(12, 6).Value = 27
$ws.Cells.Item(17, 6).Value = 1594
$ws.Cells.Item(18, 6).Value = 18
$ws.Cells.Item(20, 6).Value = 39
$ws.Cells.Item(21, 6).Value = 2664
$ws.Cells.Item(22, 6).Value = 114
$ws.Cells.Item(23, 6).Value = 1741
$ws.Cells.Item(25, 6).Value = 604
$ws.Cells.Item(26, 6).Value = 267
$ws.Cells.Item(26, 7).Value = 139
$ws.Cells.Item(28, 6).Value = 4974
$ws.Cells.Item(29, 6).Value = 63
$ws.Cells.Item(31, 6).Value = 3382
$ws.Cells.Item(32, 6).Value = 857
$ws.Cells.Item(33, 6).Value = 68
$ws.Cells.Item(35, 6).Value = 36
$ws.Cells.Item(36, 6).Value = 2408
$ws.Cells.Item(37, 6).Value = 56
$ws.Cells.Item(38, 6).Value = 363
$ws.Cells.Item(42, 6).Value = 1289
$ws.Cells.Item(45, 6).Value = 515
$ws.Cells.Item(46, 6).Value = 800
$ws.Cells.Item(47, 6).Value = 55
$ws.Cells.Item(50, 6).Value = 1410
